# Auto-generated Excel COM-interop script
# Applies per-cell numeric updates to the pricing/profit columns (H-N)
# across several worksheets, per the scheduled pricing-data refresh.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 212.5
$ws.Range("I12").Value = 212.5
$ws.Range("K12").Value = 212.5
$ws.Range("M12").Value = -42.5
# Row 18
$ws.Range("H18").Value = 440.55554
$ws.Range("I18").Value = 464.5
$ws.Range("K18").Value = 464.5
$ws.Range("M18").Value = -180.5
# Row 29
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
# Row 33
$ws.Range("H33").Value = 280.63635
$ws.Range("I33").Value = 237.44444
$ws.Range("K33").Value = 237.44444
$ws.Range("M33").Value = -8.444439999999986
# Row 40
$ws.Range("H40").Value = 1542
$ws.Range("I40").Value = 1689
$ws.Range("J40").Value = 1525.6666
$ws.Range("K40").Value = 1689
$ws.Range("L40").Value = 1525.6666
$ws.Range("M40").Value = -1514
$ws.Range("N40").Value = -1875.6666
# Row 62
$ws.Range("H62").Value = 998.5
$ws.Range("I62").Value = 998.5
$ws.Range("K62").Value = 998.5
$ws.Range("M62").Value = -374.5
# Row 65
$ws.Range("H65").Value = 998.5
$ws.Range("I65").Value = 998.5
$ws.Range("K65").Value = 4992.5
$ws.Range("M65").Value = -1872.5
# Row 70
$ws.Range("H70").Value = 3125.8572
# Row 73
$ws.Range("H73").Value = 3125.8572
# Row 103
$ws.Range("H103").Value = 1434.4
$ws.Range("I103").Value = 1407.3334
$ws.Range("K103").Value = 4222.0002
$ws.Range("M103").Value = -3636.0002
# Row 114
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
# Row 137
$ws.Range("H137").Value = 3815.6667
$ws.Range("I137").Value = 3815.6667
$ws.Range("K137").Value = 11447.0001
$ws.Range("M137").Value = -8897.000100000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 22
$ws.Range("H22").Value = 1671.6
$ws.Range("I22").Value = 1671.6
$ws.Range("K22").Value = 1671.6
$ws.Range("M22").Value = -1372.6
# Row 45
$ws.Range("H45").Value = 3364
$ws.Range("I45").Value = 2046.625
$ws.Range("K45").Value = 2046.625
$ws.Range("M45").Value = -1669.625
# Row 97
$ws.Range("H97").Value = 1111.12
$ws.Range("I97").Value = 959.5263
$ws.Range("J97").Value = 1591.1666
$ws.Range("K97").Value = 959.5263
$ws.Range("L97").Value = 1591.1666
$ws.Range("M97").Value = -463.5263
$ws.Range("N97").Value = -2583.1666
# Row 110
$ws.Range("H110").Value = 815.55554
$ws.Range("I110").Value = 667.5
$ws.Range("K110").Value = 667.5
$ws.Range("M110").Value = 1377.5
# Row 122
$ws.Range("H122").Value = 3491.5
$ws.Range("I122").Value = 4121
$ws.Range("J122").Value = 344
$ws.Range("K122").Value = 12363
$ws.Range("L122").Value = 1032
$ws.Range("M122").Value = -9913
$ws.Range("N122").Value = -5932

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 499.5
$ws.Range("I20").Value = 499.5
$ws.Range("K20").Value = 499.5
$ws.Range("M20").Value = -252.5
# Row 22
$ws.Range("H22").Value = 1074.75
$ws.Range("J22").Value = 1200
$ws.Range("L22").Value = 1200
$ws.Range("N22").Value = -1546
# Row 99
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
# Row 134
$ws.Range("H134").Value = 2725
$ws.Range("I134").Value = 2725
$ws.Range("K134").Value = 8175
$ws.Range("M134").Value = -5640
# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1604.2858
$ws.Range("I16").Value = 1046
$ws.Range("K16").Value = 1046
$ws.Range("M16").Value = -759
# Row 17
$ws.Range("H17").Value = 4003.5
$ws.Range("I17").Value = 3008
$ws.Range("J17").Value = 4999
$ws.Range("K17").Value = 3008
$ws.Range("L17").Value = 4999
$ws.Range("M17").Value = -2834
$ws.Range("N17").Value = -5347
# Row 50
$ws.Range("H50").Value = 500
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
# Row 59
$ws.Range("H59").Value = 50000
$ws.Range("J59").Value = 50000
$ws.Range("L59").Value = 50000
$ws.Range("N59").Value = -52290
# Row 62
$ws.Range("H62").Value = 4499
$ws.Range("I62").Value = 4499
$ws.Range("K62").Value = 4499
$ws.Range("M62").Value = -3875
# Row 65
$ws.Range("H65").Value = 4499
$ws.Range("I65").Value = 4499
$ws.Range("K65").Value = 22495
$ws.Range("M65").Value = -19375
# Row 68
$ws.Range("H68").Value = 90000
$ws.Range("I68").Value = 90000
$ws.Range("K68").Value = 90000
$ws.Range("M68").Value = -89251
# Row 71
$ws.Range("H71").Value = 90000
$ws.Range("I71").Value = 90000
$ws.Range("K71").Value = 270000
$ws.Range("M71").Value = -266256
# Row 107
$ws.Range("H107").Value = 978.8
$ws.Range("I107").Value = 798.3333
$ws.Range("K107").Value = 798.3333
$ws.Range("M107").Value = 1121.6667
# Row 113
$ws.Range("H113").Value = 1604.2858
$ws.Range("I113").Value = 1046
$ws.Range("K113").Value = 1046
$ws.Range("M113").Value = 1124
# Row 122
$ws.Range("H122").Value = 1004.7143
$ws.Range("I122").Value = 982.2
$ws.Range("J122").Value = 1061
$ws.Range("K122").Value = 2946.6
$ws.Range("L122").Value = 3183
$ws.Range("M122").Value = -496.6000000000004
$ws.Range("N122").Value = -8083
# Row 132
$ws.Range("H132").Value = 1901.1333
$ws.Range("I132").Value = 1672.3334
$ws.Range("J132").Value = 2244.3333
$ws.Range("K132").Value = 5017.0002
$ws.Range("L132").Value = 6732.999899999999
$ws.Range("M132").Value = -2487.0002
$ws.Range("N132").Value = -11792.9999
# Row 134
$ws.Range("H134").Value = 3834.4707
$ws.Range("I134").Value = 1562.3572
$ws.Range("J134").Value = 14437.667
$ws.Range("K134").Value = 4687.071599999999
$ws.Range("L134").Value = 43313.001
$ws.Range("M134").Value = -2152.071599999999
$ws.Range("N134").Value = -48383.001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 117
$ws.Range("H117").Value = 2771
$ws.Range("I117").Value = 1029
$ws.Range("J117").Value = 3119.4
$ws.Range("K117").Value = 3087
$ws.Range("L117").Value = 9358.200000000001
$ws.Range("M117").Value = 355
$ws.Range("N117").Value = -16242.2

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 751462.5
$ws.Range("I97").Value = 425
$ws.Range("K97").Value = 425
$ws.Range("M97").Value = 71
# Row 113
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
# Row 132
$ws.Range("H132").Value = 2985
$ws.Range("I132").Value = 2985
$ws.Range("K132").Value = 8955
$ws.Range("M132").Value = -6425

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 5874.75
$ws.Range("I7").Value = 4499.5
$ws.Range("K7").Value = 4499.5
$ws.Range("M7").Value = -4387.5
# Row 22
$ws.Range("H22").Value = 978.8
$ws.Range("J22").Value = 972.5
$ws.Range("L22").Value = 972.5
$ws.Range("N22").Value = -1562.5
# Row 27
$ws.Range("H27").Value = 978.8
$ws.Range("J27").Value = 972.5
$ws.Range("L27").Value = 972.5
$ws.Range("N27").Value = -1186.5
# Row 46
$ws.Range("H46").Value = 3223.5
$ws.Range("I46").Value = 964.6667
$ws.Range("J46").Value = 10000
$ws.Range("K46").Value = 964.6667
$ws.Range("L46").Value = 10000
$ws.Range("M46").Value = -776.6667
$ws.Range("N46").Value = -10376
# Row 68
$ws.Range("H68").Value = 2000
$ws.Range("I68").Value = 2000
$ws.Range("K68").Value = 2000
$ws.Range("M68").Value = -1251
# Row 71
$ws.Range("H71").Value = 2000
$ws.Range("I71").Value = 2000
$ws.Range("K71").Value = 10000
$ws.Range("M71").Value = -6256
# Row 122
$ws.Range("H122").Value = 2895
$ws.Range("J122").Value = 2895
$ws.Range("L122").Value = 8685
$ws.Range("N122").Value = -13585
# Row 126
$ws.Range("H126").Value = 5874.75
$ws.Range("I126").Value = 4499.5
$ws.Range("K126").Value = 13498.5
$ws.Range("M126").Value = -11028.5
# Row 141
$ws.Range("H141").Value = 30000
$ws.Range("I141").Value = 30000
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 30000
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -24820
$ws.Range("N141").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 3000
$ws.Range("I62").Value = 3000
$ws.Range("K62").Value = 3000
$ws.Range("M62").Value = -2376
# Row 65
$ws.Range("H65").Value = 3000
$ws.Range("I65").Value = 3000
$ws.Range("K65").Value = 15000
$ws.Range("M65").Value = -11880
# Row 96
$ws.Range("H96").Value = 2325
$ws.Range("I96").Value = 2325
$ws.Range("K96").Value = 2325
$ws.Range("M96").Value = -952
# Row 133
$ws.Range("H133").Value = 10714.333
$ws.Range("J133").Value = 10714.333
$ws.Range("L133").Value = 10714.333
$ws.Range("N133").Value = -20834.333

